$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New point data (elevations from ASC files under 100m are now included,
# so the whole point table was regenerated) replacing the old 10 rows
# with this new set of 5 rows.
$data = @(
    @("6167P", "289682.9416", "6190941.7817", "438.5172"),
    @("CRDX",  "294608.3090", "6199571.8370", "380.3520"),
    @("RO",    "289670.8344", "6192251.0810", "397.1787"),
    @("STN01", "289623.4902", "6192201.6504", "395.6510"),
    @("SW13",  "289565.8562", "6192191.0538", "395.7413")
)

# Drop the rows that no longer exist (old table had 10 data rows, new one has 5)
$ws.Range("A7:D11").ClearContents()

# Write the new data starting at row 2, forcing numeric-looking values to be
# stored as text (matching the original workbook's convention of storing
# every value, including numbers, as shared-string text). A leading
# apostrophe forces Excel to keep the literal text instead of coercing it
# to a number; resetting the style back to Normal afterwards keeps the
# cell format identical to the rest of the sheet.
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    for ($c = 1; $c -le 4; $c++) {
        $cell = $ws.Cells.Item($row, $c)
        $cell.Value = "'" + $data[$i][$c - 1]
    }
}
$ws.Range("A2:D6").Style = "Normal"
